$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.974.59"
$ws.Range("E2").Value = "  +12.08%  "

$ws.Range("D3").Value = "1.759.32"
$ws.Range("E3").Value = "  +8.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9954"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9927"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3849"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3661"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +20.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.238"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07732"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9939"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.484"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.104"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.93%  "

$ws.Range("D16").Value = "1.745.22"
$ws.Range("E16").Value = "  +7.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001169"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9926"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06851"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "87.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.480"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.12%  "

$ws.Range("D24").Value = "25.834.43"
$ws.Range("E24").Value = "  +11.54%  "

$ws.Range("E25").Value = "  +2.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.956"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "134.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.16%  "

$ws.Range("D30").Value = "1.934.59"
$ws.Range("E30").Value = "  +7.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.218"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +24.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.064"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +16.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.327"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "14.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +20.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.802"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08715"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.658"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06770"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.402"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02476"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2234"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.302"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6578"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9929"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6384"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.47%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.909"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.185"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07504"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.77%  "

